# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" table swaps the "Valor Mora" figures between the
# 2209 row (row 16) and the 2203 row (row 22): 2209 now carries the
# standard 84000 value, while 2203 carries the 72800 value that used
# to belong to 2209.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("F16").Value = 84000
$ws.Range("F22").Value = 72800
